# Auto-generated script applying scheduled-runner market data refresh
# to the Behemoth_Profits workbook (8 sheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1772.5454
$ws.Range("J17").Value = 1772.5454
$ws.Range("L17").Value = 5317.6362
$ws.Range("N17").Value = -5653.6362
$ws.Range("H18").Value = 1477.1666
$ws.Range("I18").Value = 1234.6471
$ws.Range("K18").Value = 1234.6471
$ws.Range("M18").Value = -950.6470999999999
$ws.Range("H32").Value = 3164.2727
$ws.Range("I32").Value = 2663.5
$ws.Range("J32").Value = 4499.6665
$ws.Range("K32").Value = 2663.5
$ws.Range("L32").Value = 4499.6665
$ws.Range("M32").Value = -2337.5
$ws.Range("N32").Value = -5151.6665
$ws.Range("H62").Value = 13109.637
$ws.Range("I62").Value = 4742
$ws.Range("K62").Value = 4742
$ws.Range("M62").Value = -4118
$ws.Range("H64").Value = 4880.75
$ws.Range("I64").Value = 4642.25
$ws.Range("J64").Value = 5000
$ws.Range("K64").Value = 4642.25
$ws.Range("L64").Value = 5000
$ws.Range("M64").Value = -4394.25
$ws.Range("N64").Value = -5496
$ws.Range("H65").Value = 13109.637
$ws.Range("I65").Value = 4742
$ws.Range("K65").Value = 23710
$ws.Range("M65").Value = -20590
$ws.Range("H67").Value = 4880.75
$ws.Range("I67").Value = 4642.25
$ws.Range("J67").Value = 5000
$ws.Range("K67").Value = 4642.25
$ws.Range("L67").Value = 5000
$ws.Range("M67").Value = -3784.25
$ws.Range("N67").Value = -6716
$ws.Range("H70").Value = 2129.625
$ws.Range("I70").Value = 2162.1
$ws.Range("J70").Value = 2075.5
$ws.Range("K70").Value = 6486.299999999999
$ws.Range("L70").Value = 6226.5
$ws.Range("M70").Value = -6216.299999999999
$ws.Range("N70").Value = -6766.5
$ws.Range("H73").Value = 2129.625
$ws.Range("I73").Value = 2162.1
$ws.Range("J73").Value = 2075.5
$ws.Range("K73").Value = 6486.299999999999
$ws.Range("L73").Value = 6226.5
$ws.Range("M73").Value = -5550.299999999999
$ws.Range("N73").Value = -8098.5
$ws.Range("H76").Value = 7233.3335
$ws.Range("I76").Value = 4400
$ws.Range("K76").Value = 4400
$ws.Range("M76").Value = -4085
$ws.Range("H79").Value = 7233.3335
$ws.Range("I79").Value = 4400
$ws.Range("K79").Value = 4400
$ws.Range("M79").Value = -3308
$ws.Range("H88").Value = 598282.1
$ws.Range("I88").Value = 41159.8
$ws.Range("J88").Value = 907794.5600000001
$ws.Range("K88").Value = 41159.8
$ws.Range("L88").Value = 907794.5600000001
$ws.Range("M88").Value = -40753.8
$ws.Range("N88").Value = -908606.5600000001
$ws.Range("H91").Value = 598282.1
$ws.Range("I91").Value = 41159.8
$ws.Range("J91").Value = 907794.5600000001
$ws.Range("K91").Value = 41159.8
$ws.Range("L91").Value = 907794.5600000001
$ws.Range("M91").Value = -39755.8
$ws.Range("N91").Value = -910602.5600000001
$ws.Range("H100").Value = 2676.875
$ws.Range("I100").Value = 1563.375
$ws.Range("K100").Value = 1563.375
$ws.Range("M100").Value = -1022.375
$ws.Range("H103").Value = 2483.4707
$ws.Range("I103").Value = 1212.5
$ws.Range("J103").Value = 3613.2222
$ws.Range("K103").Value = 3637.5
$ws.Range("L103").Value = 10839.6666
$ws.Range("M103").Value = -3051.5
$ws.Range("N103").Value = -12011.6666
$ws.Range("H104").Value = 1999
$ws.Range("J104").Value = 1999
$ws.Range("L104").Value = 5997
$ws.Range("N104").Value = -9491
$ws.Range("H116").Value = 9151.25
$ws.Range("I116").Value = 9302.5
$ws.Range("K116").Value = 9302.5
$ws.Range("M116").Value = -5860.5
$ws.Range("H132").Value = 2307.6924
$ws.Range("I132").Value = 1919.25
$ws.Range("K132").Value = 5757.75
$ws.Range("M132").Value = -3227.75
$ws.Range("H135").Value = 2841.1875
$ws.Range("I135").Value = 2675.6428
$ws.Range("J135").Value = 4000
$ws.Range("K135").Value = 24080.7852
$ws.Range("L135").Value = 36000
$ws.Range("M135").Value = -21545.7852
$ws.Range("N135").Value = -41070
$ws.Range("H138").Value = 2172.51
$ws.Range("I138").Value = 1505.3077
$ws.Range("J138").Value = 2272.2068
$ws.Range("K138").Value = 4515.9231
$ws.Range("L138").Value = 6816.6204
$ws.Range("M138").Value = 624.0769
$ws.Range("N138").Value = -17096.6204

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 858
$ws.Range("I2").Value = 551.4286
$ws.Range("K2").Value = 551.4286
$ws.Range("M2").Value = -438.4286
$ws.Range("H4").Value = 507.57144
$ws.Range("I4").Value = 530.2
$ws.Range("K4").Value = 530.2
$ws.Range("M4").Value = -414.2
$ws.Range("H32").Value = 25009156
$ws.Range("I32").Value = 27781838
$ws.Range("K32").Value = 27781838
$ws.Range("M32").Value = -27781551
$ws.Range("H33").Value = 15000
$ws.Range("I33").Value = 15000
$ws.Range("K33").Value = 15000
$ws.Range("M33").Value = -14671
$ws.Range("H45").Value = 50002970
$ws.Range("I45").Value = 100001680
$ws.Range("K45").Value = 100001680
$ws.Range("M45").Value = -100001303
$ws.Range("H74").Value = 14450837
$ws.Range("I74").Value = 31254612
$ws.Range("K74").Value = 31254612
$ws.Range("M74").Value = -31253738
$ws.Range("H77").Value = 14450837
$ws.Range("I77").Value = 31254612
$ws.Range("K77").Value = 156273060
$ws.Range("M77").Value = -156268692
$ws.Range("H88").Value = 2876.276
$ws.Range("I88").Value = 2348
$ws.Range("J88").Value = 3014.087
$ws.Range("K88").Value = 2348
$ws.Range("L88").Value = 3014.087
$ws.Range("M88").Value = -1942
$ws.Range("N88").Value = -3826.087
$ws.Range("H91").Value = 2876.276
$ws.Range("I91").Value = 2348
$ws.Range("J91").Value = 3014.087
$ws.Range("K91").Value = 2348
$ws.Range("L91").Value = 3014.087
$ws.Range("M91").Value = -944
$ws.Range("N91").Value = -5822.087
$ws.Range("H102").Value = 27046.166
$ws.Range("I102").Value = 27046.166
$ws.Range("K102").Value = 27046.166
$ws.Range("M102").Value = -25424.166
$ws.Range("H116").Value = 858
$ws.Range("I116").Value = 551.4286
$ws.Range("K116").Value = 551.4286
$ws.Range("M116").Value = 1742.5714
$ws.Range("H122").Value = 2969
$ws.Range("I122").Value = 2717.1428
$ws.Range("J122").Value = 3321.6
$ws.Range("K122").Value = 8151.428400000001
$ws.Range("L122").Value = 9964.799999999999
$ws.Range("M122").Value = -5701.428400000001
$ws.Range("N122").Value = -14864.8

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 858
$ws.Range("I3").Value = 551.4286
$ws.Range("K3").Value = 551.4286
$ws.Range("M3").Value = -437.4286
$ws.Range("H20").Value = 9966.333000000001
$ws.Range("I20").Value = 9966.333000000001
$ws.Range("K20").Value = 9966.333000000001
$ws.Range("M20").Value = -9719.333000000001
$ws.Range("H40").Value = 177963
$ws.Range("J40").Value = 177963
$ws.Range("L40").Value = 177963
$ws.Range("N40").Value = -178493
$ws.Range("H94").Value = 2424.5
$ws.Range("I94").Value = 3665.3333
$ws.Range("J94").Value = 1892.7142
$ws.Range("K94").Value = 3665.3333
$ws.Range("L94").Value = 1892.7142
$ws.Range("M94").Value = -3214.3333
$ws.Range("N94").Value = -2794.7142
$ws.Range("H96").Value = 62161.668
$ws.Range("I96").Value = 14904.5
$ws.Range("J96").Value = 109418.836
$ws.Range("K96").Value = 14904.5
$ws.Range("L96").Value = 109418.836
$ws.Range("M96").Value = -12158.5
$ws.Range("N96").Value = -114910.836
$ws.Range("H97").Value = 35497.25
$ws.Range("I97").Value = 21999.5
$ws.Range("K97").Value = 21999.5
$ws.Range("M97").Value = -21008.5
$ws.Range("H99").Value = 3015.487
$ws.Range("I99").Value = 2356.28
$ws.Range("K99").Value = 2356.28
$ws.Range("M99").Value = -858.2800000000002

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 69539.336
$ws.Range("J18").Value = 69539.336
$ws.Range("L18").Value = 69539.336
$ws.Range("N18").Value = -69999.336
$ws.Range("H22").Value = 468.66666
$ws.Range("I22").Value = 399.6
$ws.Range("J22").Value = 555
$ws.Range("K22").Value = 399.6
$ws.Range("L22").Value = 555
$ws.Range("M22").Value = -49.60000000000002
$ws.Range("N22").Value = -1255
$ws.Range("H33").Value = 3255.375
$ws.Range("I33").Value = 3149
$ws.Range("J33").Value = 4000
$ws.Range("K33").Value = 3149
$ws.Range("L33").Value = 4000
$ws.Range("M33").Value = -2770
$ws.Range("N33").Value = -4758
$ws.Range("H35").Value = 2933.3333
$ws.Range("I35").Value = 1900
$ws.Range("J35").Value = 5000
$ws.Range("K35").Value = 1900
$ws.Range("L35").Value = 5000
$ws.Range("M35").Value = -1606
$ws.Range("N35").Value = -5588
$ws.Range("H36").Value = 2225
$ws.Range("J36").Value = 4000
$ws.Range("L36").Value = 4000
$ws.Range("N36").Value = -4776
$ws.Range("H37").Value = 3999.5
$ws.Range("J37").Value = 3999.5
$ws.Range("L37").Value = 3999.5
$ws.Range("N37").Value = -4213.5
$ws.Range("H40").Value = 2225
$ws.Range("J40").Value = 4000
$ws.Range("L40").Value = 4000
$ws.Range("N40").Value = -4320
$ws.Range("H62").Value = 2506480.8
$ws.Range("I62").Value = 2506480.8
$ws.Range("K62").Value = 2506480.8
$ws.Range("M62").Value = -2505856.8
$ws.Range("H65").Value = 2506480.8
$ws.Range("I65").Value = 2506480.8
$ws.Range("K65").Value = 12532404
$ws.Range("M65").Value = -12529284
$ws.Range("H94").Value = 4345
$ws.Range("I94").Value = 2324.3333
$ws.Range("K94").Value = 2324.3333
$ws.Range("M94").Value = -1873.3333
$ws.Range("H99").Value = 4807
$ws.Range("I99").Value = 4807
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 4807
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -3309
$ws.Range("N99").ClearContents()
$ws.Range("H122").Value = 2577.3333
$ws.Range("I122").Value = 2712.8
$ws.Range("K122").Value = 8138.400000000001
$ws.Range("M122").Value = -5688.400000000001
$ws.Range("H126").Value = 4807
$ws.Range("I126").Value = 4807
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 14421
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -11951
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 1772.9
$ws.Range("I132").Value = 1772.9
$ws.Range("K132").Value = 5318.700000000001
$ws.Range("M132").Value = -2788.700000000001
$ws.Range("H134").Value = 5694.077
$ws.Range("I134").Value = 2689
$ws.Range("J134").Value = 8269.857
$ws.Range("K134").Value = 8067
$ws.Range("L134").Value = 24809.571
$ws.Range("M134").Value = -5532
$ws.Range("N134").Value = -29879.571
$ws.Range("H141").Value = 229372.1
$ws.Range("J141").Value = 242409.3
$ws.Range("L141").Value = 242409.3
$ws.Range("N141").Value = -252769.3

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3824.8333
$ws.Range("J68").Value = 4390
$ws.Range("L68").Value = 13170
$ws.Range("N68").Value = -14792
$ws.Range("H71").Value = 3824.8333
$ws.Range("J71").Value = 4390
$ws.Range("L71").Value = 39510
$ws.Range("N71").Value = -47622
$ws.Range("H80").Value = 4303.25
$ws.Range("J80").Value = 4294.273
$ws.Range("L80").Value = 12882.819
$ws.Range("N80").Value = -14754.819
$ws.Range("H83").Value = 4303.25
$ws.Range("J83").Value = 4294.273
$ws.Range("L83").Value = 38648.457
$ws.Range("N83").Value = -48008.457
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H108").Value = 0
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("M108").ClearContents()
$ws.Range("N108").ClearContents()
$ws.Range("H114").Value = 2000
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H131").Value = 5717.089
$ws.Range("J131").Value = 5366.237
$ws.Range("L131").Value = 16098.711
$ws.Range("N131").Value = -26178.711
$ws.Range("H132").Value = 1858
$ws.Range("I132").Value = 1802.4
$ws.Range("J132").Value = 1888.8889
$ws.Range("K132").Value = 16221.6
$ws.Range("L132").Value = 17000.0001
$ws.Range("M132").Value = -13691.6
$ws.Range("N132").Value = -22060.0001
$ws.Range("H137").Value = 5317.5
$ws.Range("J137").Value = 5247
$ws.Range("L137").Value = 15741
$ws.Range("N137").Value = -25941

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 100983.664
$ws.Range("I5").Value = 146475
$ws.Range("K5").Value = 146475
$ws.Range("M5").Value = -146363
$ws.Range("H80").Value = 3024.842
$ws.Range("I80").Value = 3049.5833
$ws.Range("J80").Value = 2982.4285
$ws.Range("K80").Value = 3049.5833
$ws.Range("L80").Value = 2982.4285
$ws.Range("M80").Value = -2051.5833
$ws.Range("N80").Value = -4978.4285
$ws.Range("H83").Value = 3024.842
$ws.Range("I83").Value = 3049.5833
$ws.Range("J83").Value = 2982.4285
$ws.Range("K83").Value = 15247.9165
$ws.Range("L83").Value = 14912.1425
$ws.Range("M83").Value = -10255.9165
$ws.Range("N83").Value = -24896.1425
$ws.Range("H97").Value = 1973.5454
$ws.Range("I97").Value = 3008.1667
$ws.Range("J97").Value = 732
$ws.Range("K97").Value = 3008.1667
$ws.Range("L97").Value = 732
$ws.Range("M97").Value = -2512.1667
$ws.Range("N97").Value = -1724
$ws.Range("H99").Value = 24872.2
$ws.Range("I99").Value = 14737.25
$ws.Range("J99").Value = 65412
$ws.Range("K99").Value = 14737.25
$ws.Range("L99").Value = 65412
$ws.Range("M99").Value = -12491.25
$ws.Range("N99").Value = -69904
$ws.Range("H107").Value = 389.66666
$ws.Range("I107").Value = 353.2
$ws.Range("J107").Value = 572
$ws.Range("K107").Value = 353.2
$ws.Range("L107").Value = 572
$ws.Range("M107").Value = 1566.8
$ws.Range("N107").Value = -4412
$ws.Range("H108").Value = 99819.25
$ws.Range("J108").Value = 99819.25
$ws.Range("L108").Value = 99819.25
$ws.Range("N108").Value = -107499.25
$ws.Range("H113").Value = 3147.2273
$ws.Range("I113").Value = 2143.9092
$ws.Range("K113").Value = 2143.9092
$ws.Range("M113").Value = 26.09079999999994
$ws.Range("H122").Value = 1894.9166
$ws.Range("I122").Value = 1641.8572
$ws.Range("K122").Value = 4925.571599999999
$ws.Range("M122").Value = -2475.571599999999
$ws.Range("H126").Value = 3248.7222
$ws.Range("I126").Value = 2477.6
$ws.Range("J126").Value = 3545.3076
$ws.Range("K126").Value = 7432.799999999999
$ws.Range("L126").Value = 10635.9228
$ws.Range("M126").Value = -4962.799999999999
$ws.Range("N126").Value = -15575.9228
$ws.Range("H132").Value = 45465164
$ws.Range("I132").Value = 52637790
$ws.Range("K132").Value = 157913370
$ws.Range("M132").Value = -157910840
$ws.Range("H134").Value = 72499.75
$ws.Range("J134").Value = 72499.75
$ws.Range("L134").Value = 217499.25
$ws.Range("N134").Value = -222569.25
$ws.Range("H136").Value = 27999
$ws.Range("J136").Value = 27999
$ws.Range("L136").Value = 83997
$ws.Range("N136").Value = -89097

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1208.25
$ws.Range("I22").Value = 1039.6
$ws.Range("J22").Value = 1328.7142
$ws.Range("K22").Value = 1039.6
$ws.Range("L22").Value = 1328.7142
$ws.Range("M22").Value = -744.5999999999999
$ws.Range("N22").Value = -1918.7142
$ws.Range("H27").Value = 1208.25
$ws.Range("I27").Value = 1039.6
$ws.Range("J27").Value = 1328.7142
$ws.Range("K27").Value = 1039.6
$ws.Range("L27").Value = 1328.7142
$ws.Range("M27").Value = -932.5999999999999
$ws.Range("N27").Value = -1542.7142
$ws.Range("H38").Value = 54999
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("H40").Value = 3224.4
$ws.Range("I40").Value = 2659.1072
$ws.Range("K40").Value = 2659.1072
$ws.Range("M40").Value = -2523.1072
$ws.Range("H46").Value = 5194.8
$ws.Range("I46").Value = 2092.4285
$ws.Range("K46").Value = 2092.4285
$ws.Range("M46").Value = -1904.4285
$ws.Range("H87").Value = 70000
$ws.Range("J87").Value = 100000
$ws.Range("L87").Value = 100000
$ws.Range("N87").Value = -102246
$ws.Range("H90").Value = 70000
$ws.Range("J90").Value = 100000
$ws.Range("L90").Value = 300000
$ws.Range("N90").Value = -311232
$ws.Range("H93").Value = 62501290
$ws.Range("I93").Value = 100001144
$ws.Range("K93").Value = 100001144
$ws.Range("M93").Value = -99999896
$ws.Range("H122").Value = 4510.1577
$ws.Range("I122").Value = 4102.8
$ws.Range("K122").Value = 12308.4
$ws.Range("M122").Value = -9858.400000000001
$ws.Range("H136").Value = 92205.69
$ws.Range("I136").Value = 15811.714
$ws.Range("K136").Value = 47435.142
$ws.Range("M136").Value = -44885.142

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 29666.666
$ws.Range("J25").Value = 29666.666
$ws.Range("L25").Value = 29666.666
$ws.Range("N25").Value = -30252.666
$ws.Range("H27").Value = 126000
$ws.Range("J27").Value = 126000
$ws.Range("L27").Value = 126000
$ws.Range("N27").Value = -126138
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H122").Value = 2843.2144
$ws.Range("I122").Value = 3140.6
$ws.Range("K122").Value = 9421.799999999999
$ws.Range("M122").Value = -6971.799999999999
$ws.Range("H123").Value = 54569.855
$ws.Range("J123").Value = 54569.855
$ws.Range("L123").Value = 54569.855
$ws.Range("N123").Value = -64369.855
$ws.Range("H132").Value = 1710.8914
$ws.Range("I132").Value = 1282.7142
$ws.Range("K132").Value = 3848.1426
$ws.Range("M132").Value = -1318.1426

Write-Host "Applied 461 cell updates and 8 cell clears across 8 sheets."
